# Adds a new "Q" column (one more day of cumulative data) to the daily
# "Fallecidos Min Ciencias acumulado" tracker, mirroring column P's
# layout/formatting, and appends a new trailing row (91) for date 44002.
#
# Also widens the O1/P1 running-total formulas from a fixed end row (88)
# to an open-ended 999 so future appended rows keep summing automatically,
# and adds the same open-ended SUM formula for the new Q1 total cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Wire up the running-total formulas *before* touching any of the
#    underlying data cells, so the engine's recalculation graph picks up
#    every later edit to O3:O999 / P3:P999 / Q3:Q999 (setting them only
#    after the data is in place can leave the cached <v> one row stale).
# ---------------------------------------------------------------------
$ws.Range("O1").Formula = "=SUM(O3:O999)"
$ws.Range("P1").Formula = "=SUM(P3:P999)"
$ws.Range("Q1").Formula = "=SUM(Q3:Q999)"

# ---------------------------------------------------------------------
# 2) Clone column P's formatting into column Q (header row through the
#    last existing data row) so every Q cell picks up the same style
#    index as its P neighbour (totals row, date row, data rows, ...).
# ---------------------------------------------------------------------
$ws.Range("P1:P90").Copy() | Out-Null
$ws.Range("Q1:Q90").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Clone column A's formatting (date style) down into the new row 91,
#    and column P's data-cell formatting into the new Q91 cell.
# ---------------------------------------------------------------------
$ws.Range("A90").Copy() | Out-Null
$ws.Range("A91").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P90").Copy() | Out-Null
$ws.Range("Q91").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Populate the new date cells.
# ---------------------------------------------------------------------
$ws.Range("Q2").Value = 44004
$ws.Range("A91").Value = 44002

# ---------------------------------------------------------------------
# 5) Populate the Q data column (rows 3-90), one extra day's worth of
#    cumulative counts - mostly identical to column P, save for a
#    handful of rows where the count was revised upward.
# ---------------------------------------------------------------------
$qData = @(1,1,1,1,1,1,3,4,1,4,6,8,1,5,4,10,9,6,7,3,9,4,8,6,5,6,5,10,8,5,7,6,7,9,11,7,6,8,10,7,19,12,11,14,30,16,22,30,19,26,38,35,32,49,46,63,61,65,62,70,86,95,86,106,118,113,121,126,113,142,133,149,144,148,155,144,135,168,153,158,152,153,147,134,153,128,80,20)
$qArr = New-Object 'object[,]' $qData.Length,1
for ($i = 0; $i -lt $qData.Length; $i++) {
    $qArr[$i,0] = $qData[$i]
}
$ws.Range("Q3:Q90").Value2 = $qArr

# New trailing row for the next day, only the count is known so far.
$ws.Range("Q91").Value = 1

# ---------------------------------------------------------------------
# 6) Restore view state (best effort - scroll position / selection).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("S8").Select()

Write-Host "edit complete"
